$d = $word.ActiveDocument

# Change 1: merge runs in "Objectifs" paragraph (no text change, just run merge)
$d.Content.Find.Execute(" : dans ce chapitre, vous apprendrez à …", $true, $false, $false, $false, $false, $true, 1, $false, " : dans ce chapitre, vous apprendrez à …", 2) | Out-Null

# Change 2: merge runs in "Sommaire" paragraph (no text change, just run merge)
$d.Content.Find.Execute("Introduction – Premier problème…", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction – Premier problème…", 2) | Out-Null

# Change 3: remove the _GoBack bookmark from the "Remarque : ce style..." paragraph
$d.Bookmarks.Item("_GoBack").Delete()

# Change 4: merge + simplify wording "en bleu clair, en gras et en italique" -> "en bleu clair, gras et italique"
$d.Content.Find.Execute("apparaissent, hors listing, en bleu clair, en gras et en italique : ", $true, $false, $false, $false, $false, $true, 1, $false, "apparaissent, hors listing, en bleu clair, gras et italique : ", 2) | Out-Null

# Change 5: simplify wording "en gras et en italiques" -> "en gras et italique", then
# split the run and re-insert the _GoBack bookmark right after "italique"
$d.Content.Find.Execute("apparaissent en gras et en italiques", $true, $false, $false, $false, $false, $true, 1, $false, "apparaissent en gras et italique", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("apparaissent en gras et italique", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)
